# Apply the January-2024 data refresh described in the commit.
#
# Sheet "部门情况202401" (department overview): several department rows'
# outstanding-balance / account-count / rate figures were revised. These
# cells hold numeric-looking text (they were authored as text, not numbers),
# so we use a leading apostrophe to force text entry and avoid Excel's
# automatic number coercion.
#
# Sheet "对公产品台账202401" (corporate product ledger): the blank-name
# summary row and the "一般贷款" row at the top were removed, and the
# trailing "线下卖方产品" row was removed as well, shifting everything
# else up. The trailing "Total" row was recalculated to reflect the sum
# of the remaining detail rows.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 部门情况202401 ---------------------------------------------
$ws1 = $wb.Worksheets.Item("部门情况202401")

$ws1.Range("B2").Value = "'152526.63"
$ws1.Range("C2").Value = "'258.00"

$ws1.Range("B3").Value = "'141173.17"
$ws1.Range("C3").Value = "'32.00"
$ws1.Range("E3").Value = "'5.77"

$ws1.Range("B4").Value = "'143196.64"
$ws1.Range("C4").Value = "'115.00"

$ws1.Range("B5").Value = "'60898.34"
$ws1.Range("C5").Value = "'13.00"

$ws1.Range("B6").Value = "'8500.00"
$ws1.Range("C6").Value = "'2.00"

# --- Sheet 4: 对公产品台账202401 -----------------------------------------
$ws4 = $wb.Worksheets.Item("对公产品台账202401")

# Remove the blank-name summary row (old row 2) and the "一般贷款" row
# (old row 3). Deleting twice at index 2 removes both, shifting all the
# detail rows up by two.
$ws4.Rows(2).Delete()
$ws4.Rows(2).Delete()

# Remove the trailing "线下卖方产品" row, now at row 15 (was row 17).
$ws4.Rows(15).Delete()

# Recalculate the "Total" row (now row 15, was row 18) over the remaining
# detail rows.
$ws4.Range("C15").Value = 227729.11
$ws4.Range("D15").Value = 399
$ws4.Range("E15").Value = 480063
$ws4.Range("F15").Value = 29392.56
$ws4.Range("G15").Value = 133
$ws4.Range("H15").Value = 141320.63
$ws4.Range("I15").Value = 295
$ws4.Range("J15").Value = 5001.79
$ws4.Range("K15").Value = 110.16
